# The workbook contained the below-50%-MIT-score guide RNA rows (1-10)
# duplicated a second time (rows 11-20), followed by three extra rows
# (21-23 here originally at 23-25) that actually belonged to an "above 50"
# / "no upstream" result set that leaked into this file. Per the commit
# message ("above 50 no upstream 1 file") those trailing rows need to go,
# and the sheet should end up holding just the original 10 unique rows
# repeated twice (20 rows total), matching the new dimension A1:S20 and
# the shrunken shared-strings table (uniqueCount 112 -> 80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing rows (21-25) that don't belong in this sheet.
$ws.Rows("21:25").Delete()

# Rows 11-20 need to hold the same data as rows 1-10 (the duplicate block
# that the upstream generator produces). Copy/paste the first block over
# the second to realign it.
$ws.Range("A1:S10").Copy()
$ws.Range("A11").PasteSpecial()

Write-Output "done"
